$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the missing "F" column value (5) for rows 31-35, matching the
# formatting already used by the other cells in column F (centered both
# horizontally and vertically).
foreach ($r in 31..35) {
    $dst = $ws.Cells.Item($r, 6)       # F(r)
    $dst.Value = 5
    $dst.HorizontalAlignment = -4108   # xlCenter
    $dst.VerticalAlignment = -4108     # xlCenter
}

# Update the saved view state of the sheet: scroll position and selection.
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 3
$ws.Range("D34").Select() | Out-Null
